$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.869.57"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "'1.860.27"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'304.78"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.5046"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").Value = "'0.3644"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Value = "'0.07185"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'20.69"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'1.866.80"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'0.07489"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'92.35"
$ws.Range("D15").Value = "'5.225"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'26.909.89"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'5.029"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'2.096.39"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "'10.37"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "'6.401"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").Value = "'147.81"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").Value = "'17.86"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'2.063"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "'4.683"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "'4.671"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'0.09258"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "'0.05091"
$ws.Range("D34").Value = "'0.7464"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "'2.953"
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("D36").Value = "'1.150"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'3.283"
$ws.Range("E37").Value = "  +7.79%  "
$ws.Range("D38").Value = "'0.02000"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "'2.499"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'0.5502"
$ws.Range("E40").Value = "  +3.84%  "
$ws.Range("D41").Value = "'1.070"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").Value = "'118.56"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").Value = "'6.491"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'8.476"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'0.4677"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "'0.9999"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'10.01"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "'1.566"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'36.97"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "'63.08"
$ws.Range("E51").Value = "  -2.15%  "
